{"js": "// Update the date line and the twenty-five \"two-digit \u00d7 two-digit\"\n// answer cells to the values from the next day's worksheet.\nconst replacements = [\n  [\"2025-04-25 Friday\", \"2025-04-26 Saturday\"],\n  [\"65\u00d714=910\", \"92\u00d743=3956\"],\n  [\"83\u00d795=7885\", \"53\u00d791=4823\"],\n  [\"69\u00d757=3933\", \"75\u00d720=1500\"],\n  [\"67\u00d726=1742\", \"82\u00d739=3198\"],\n  [\"72\u00d741=2952\", \"32\u00d764=2048\"],\n  [\"60\u00d741=2460\", \"71\u00d785=6035\"],\n  [\"33\u00d737=1221\", \"93\u00d743=3999\"],\n  [\"92\u00d772=6624\", \"71\u00d756=3976\"],\n  [\"71\u00d716=1136\", \"53\u00d731=1643\"],\n  [\"65\u00d772=4680\", \"20\u00d782=1640\"],\n  [\"73\u00d746=3358\", \"45\u00d719=855\"],\n  [\"62\u00d761=3782\", \"31\u00d797=3007\"],\n  [\"85\u00d718=1530\", \"47\u00d775=3525\"],\n  [\"94\u00d791=8554\", \"66\u00d712=792\"],\n  [\"49\u00d799=4851\", \"93\u00d724=2232\"],\n  [\"91\u00d737=3367\", \"65\u00d762=4030\"],\n  [\"73\u00d732=2336\", \"63\u00d789=5607\"],\n  [\"46\u00d783=3818\", \"51\u00d769=3519\"],\n  [\"79\u00d725=1975\", \"54\u00d779=4266\"],\n  [\"68\u00d716=1088\", \"60\u00d795=5700\"],\n  [\"32\u00d751=1632\", \"65\u00d766=4290\"],\n  [\"18\u00d767=1206\", \"57\u00d747=2679\"],\n  [\"25\u00d752=1300\", \"86\u00d756=4816\"],\n  [\"84\u00d711=924\", \"54\u00d732=1728\"],\n  [\"30\u00d732=960\", \"66\u00d719=1254\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"two-digit \u00d7 two-digit\"\n# answer cells to the values from the next day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-25 Friday\", \"2025-04-26 Saturday\"),\n    @(\"65\u00d714=910\", \"92\u00d743=3956\"),\n    @(\"83\u00d795=7885\", \"53\u00d791=4823\"),\n    @(\"69\u00d757=3933\", \"75\u00d720=1500\"),\n    @(\"67\u00d726=1742\", \"82\u00d739=3198\"),\n    @(\"72\u00d741=2952\", \"32\u00d764=2048\"),\n    @(\"60\u00d741=2460\", \"71\u00d785=6035\"),\n    @(\"33\u00d737=1221\", \"93\u00d743=3999\"),\n    @(\"92\u00d772=6624\", \"71\u00d756=3976\"),\n    @(\"71\u00d716=1136\", \"53\u00d731=1643\"),\n    @(\"65\u00d772=4680\", \"20\u00d782=1640\"),\n    @(\"73\u00d746=3358\", \"45\u00d719=855\"),\n    @(\"62\u00d761=3782\", \"31\u00d797=3007\"),\n    @(\"85\u00d718=1530\", \"47\u00d775=3525\"),\n    @(\"94\u00d791=8554\", \"66\u00d712=792\"),\n    @(\"49\u00d799=4851\", \"93\u00d724=2232\"),\n    @(\"91\u00d737=3367\", \"65\u00d762=4030\"),\n    @(\"73\u00d732=2336\", \"63\u00d789=5607\"),\n    @(\"46\u00d783=3818\", \"51\u00d769=3519\"),\n    @(\"79\u00d725=1975\", \"54\u00d779=4266\"),\n    @(\"68\u00d716=1088\", \"60\u00d795=5700\"),\n    @(\"32\u00d751=1632\", \"65\u00d766=4290\"),\n    @(\"18\u00d767=1206\", \"57\u00d747=2679\"),\n    @(\"25\u00d752=1300\", \"86\u00d756=4816\"),\n    @(\"84\u00d711=924\", \"54\u00d732=1728\"),\n    @(\"30\u00d732=960\", \"66\u00d719=1254\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
